$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 19.62974171293406
$ws.Range("C2").Value = 11.90515418833878
$ws.Range("D2").Value = 5.273071938883573
$ws.Range("E2").Value = 10.73751006893334
$ws.Range("F2").Value = 58.40536167705619
$ws.Range("J2").Value = 10.20398338144547
$ws.Range("M2").Value = 19.38323440442752
$ws.Range("B3").Value = 19.39771230712965
$ws.Range("C3").Value = 11.71739948170023
$ws.Range("D3").Value = 5.216086276269363
$ws.Range("E3").Value = 10.77371095808514
$ws.Range("F3").Value = 57.47542236806399
$ws.Range("J3").Value = 10.20837293539281
$ws.Range("M3").Value = 19.39223646570822
$ws.Range("B4").Value = 19.26257364858362
$ws.Range("C4").Value = 11.60724474836054
$ws.Range("D4").Value = 5.18142274897335
$ws.Range("E4").Value = 10.79727827854841
$ws.Range("F4").Value = 56.90513245647951
$ws.Range("J4").Value = 10.21193096383938
$ws.Range("M4").Value = 19.40398803934888
$ws.Range("B5").Value = 19.2094193426264
$ws.Range("C5").Value = 11.56371910496066
$ws.Range("D5").Value = 5.167386590411171
$ws.Range("E5").Value = 10.80722020616806
$ws.Range("F5").Value = 56.67312497238476
$ws.Range("J5").Value = 10.21359756955098
$ws.Range("M5").Value = 19.41033945220041
$ws.Range("B6").Value = 19.20071085665105
$ws.Range("C6").Value = 11.55657613311798
$ws.Range("D6").Value = 5.16506151776362
$ws.Range("E6").Value = 10.80889150686918
$ws.Range("F6").Value = 56.63463015741048
$ws.Range("J6").Value = 10.21388738643805
$ws.Range("M6").Value = 19.41148837811568
$ws.Range("B7").Value = 19.2618489411463
$ws.Range("C7").Value = 11.60665213131108
$ws.Range("D7").Value = 5.181233080257132
$ws.Range("E7").Value = 10.79741098856525
$ws.Range("F7").Value = 56.90200165818125
$ws.Range("J7").Value = 10.21195256331952
$ws.Range("M7").Value = 19.40406737413605
$ws.Range("B8").Value = 19.54826225844387
$ws.Range("C8").Value = 11.83939763723883
$ws.Range("D8").Value = 5.253360266723776
$ws.Range("E8").Value = 10.74971474115319
$ws.Range("F8").Value = 58.08470302662855
$ws.Range("J8").Value = 10.20531765310042
$ws.Range("M8").Value = 19.38504520643841
$ws.Range("B9").Value = 20.16463181465972
$ws.Range("C9").Value = 12.33308491229806
$ws.Range("D9").Value = 5.397043490282215
$ws.Range("E9").Value = 10.66676207137632
$ws.Range("F9").Value = 60.40010644014476
$ws.Range("J9").Value = 10.19916581532182
$ws.Range("M9").Value = 19.39720184337372
$ws.Range("B10").Value = 20.64598627138486
$ws.Range("C10").Value = 12.71381497132894
$ws.Range("D10").Value = 5.503501399746922
$ws.Range("E10").Value = 10.61219537240913
$ws.Range("F10").Value = 62.08648415599771
$ws.Range("J10").Value = 10.19884609355905
$ws.Range("M10").Value = 19.43631965603128
$ws.Range("B11").Value = 20.87008020179767
$ws.Range("C11").Value = 12.88992159134436
$ws.Range("D11").Value = 5.552022358176979
$ws.Range("E11").Value = 10.58874159962451
$ws.Range("F11").Value = 62.84791786180845
$ws.Range("J11").Value = 10.19961630304813
$ws.Range("M11").Value = 19.46065497148264
$ws.Range("B12").Value = 20.95558421132884
$ws.Range("C12").Value = 12.95694422062285
$ws.Range("D12").Value = 5.570400427781046
$ws.Range("E12").Value = 10.58005594899956
$ws.Range("F12").Value = 63.13522761613369
$ws.Range("J12").Value = 10.20003986407435
$ws.Range("M12").Value = 19.47080745234948
$ws.Range("B13").Value = 20.93714219377313
$ws.Range("C13").Value = 12.94249611465856
$ws.Range("D13").Value = 5.566442338025587
$ws.Range("E13").Value = 10.58191786851923
$ws.Range("F13").Value = 63.07339929687659
$ws.Range("J13").Value = 10.19994277294549
$ws.Range("M13").Value = 19.46857931473991
$ws.Range("B14").Value = 20.87710227862295
$ws.Range("C14").Value = 12.8954293439448
$ws.Range("D14").Value = 5.553534300203652
$ws.Range("E14").Value = 10.58802310816461
$ws.Range("F14").Value = 62.87157678723252
$ws.Range("J14").Value = 10.19964850496848
$ws.Range("M14").Value = 19.46147146159649
$ws.Range("B15").Value = 20.84040729711613
$ws.Range("C15").Value = 12.86664074978993
$ws.Range("D15").Value = 5.545628036081151
$ws.Range("E15").Value = 10.59178820925142
$ws.Range("F15").Value = 62.74781467270226
$ws.Range("J15").Value = 10.1994854409809
$ws.Range("M15").Value = 19.45723963413024
$ws.Range("B16").Value = 20.6314364137041
$ws.Range("C16").Value = 12.70235736237117
$ws.Range("D16").Value = 5.500331529110782
$ws.Range("E16").Value = 10.61375558635296
$ws.Range("F16").Value = 62.03658987126207
$ws.Range("J16").Value = 10.19881420395944
$ws.Range("M16").Value = 19.43486064130136
$ws.Range("B17").Value = 20.50448612400438
$ws.Range("C17").Value = 12.60226018931331
$ws.Range("D17").Value = 5.4725608665301
$ws.Range("E17").Value = 10.62758170934053
$ws.Range("F17").Value = 61.59866427041031
$ws.Range("J17").Value = 10.19863711269799
$ws.Range("M17").Value = 19.42280503746831
$ws.Range("B18").Value = 20.4319554509717
$ws.Range("C18").Value = 12.54496581472105
$ws.Range("D18").Value = 5.456597189551628
$ws.Range("E18").Value = 10.63566303916464
$ws.Range("F18").Value = 61.34625577150141
$ws.Range("J18").Value = 10.1986214277126
$ws.Range("M18").Value = 19.41648685714537
$ws.Range("B19").Value = 20.40748439909514
$ws.Range("C19").Value = 12.52561740052087
$ws.Range("D19").Value = 5.451194032742896
$ws.Range("E19").Value = 10.63842141136868
$ws.Range("F19").Value = 61.26071083113153
$ws.Range("J19").Value = 10.198630910023
$ws.Range("M19").Value = 19.41445348995023
$ws.Range("B20").Value = 20.51795038483312
$ws.Range("C20").Value = 12.61288745439388
$ws.Range("D20").Value = 5.475516207091395
$ws.Range("E20").Value = 10.62609656101029
$ws.Range("F20").Value = 61.64533809227889
$ws.Range("J20").Value = 10.19864704403269
$ws.Range("M20").Value = 19.42402465561613
$ws.Range("B21").Value = 20.89472069910257
$ws.Range("C21").Value = 12.90924556362604
$ws.Range("D21").Value = 5.557325661306973
$ws.Range("E21").Value = 10.58622454617695
$ws.Range("F21").Value = 62.93088648632042
$ws.Range("J21").Value = 10.19973135717284
$ws.Range("M21").Value = 19.46353380604293
$ws.Range("B22").Value = 21.14467611788627
$ws.Range("C22").Value = 13.1048474128817
$ws.Range("D22").Value = 5.610813356394927
$ws.Range("E22").Value = 10.56130657257293
$ws.Range("F22").Value = 63.76498006604643
$ws.Range("J22").Value = 10.20120892047465
$ws.Range("M22").Value = 19.49481587292935
$ws.Range("B23").Value = 21.01096071218067
$ws.Range("C23").Value = 13.0003029076296
$ws.Range("D23").Value = 5.582267092636163
$ws.Range("E23").Value = 10.57450173896866
$ws.Range("F23").Value = 63.32043217726865
$ws.Range("J23").Value = 10.20034989225893
$ws.Range("M23").Value = 19.47762177717762
$ws.Range("B24").Value = 20.51186176910112
$ws.Range("C24").Value = 12.6080820739912
$ws.Range("D24").Value = 5.474180090558031
$ws.Range("E24").Value = 10.62676758375858
$ws.Range("F24").Value = 61.62423883264124
$ws.Range("J24").Value = 10.19864228582705
$ws.Range("M24").Value = 19.42347135744638
$ws.Range("B25").Value = 19.99252008330682
$ws.Range("C25").Value = 12.19602102589518
$ws.Range("D25").Value = 5.357981458067443
$ws.Range("E25").Value = 10.68807771864415
$ws.Range("F25").Value = 59.77543653153571
$ws.Range("J25").Value = 10.20009373353879
$ws.Range("M25").Value = 19.38861031543715
